$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3609.182
$ws.Range("I40").Value = 2499.75
$ws.Range("K40").Value = 2499.75
$ws.Range("M40").Value = -2324.75
$ws.Range("H55").Value = 863.7143
$ws.Range("J55").Value = 1079.8
$ws.Range("L55").Value = 1079.8
$ws.Range("N55").Value = -1507.8
$ws.Range("H95").Value = 24999.5
$ws.Range("J95").Value = 24999.5
$ws.Range("L95").Value = 24999.5
$ws.Range("N95").Value = -30491.5
$ws.Range("H107").Value = 1448.7727
$ws.Range("I107").Value = 1537.2142
$ws.Range("J107").Value = 1294
$ws.Range("K107").Value = 1537.2142
$ws.Range("L107").Value = 1294
$ws.Range("M107").Value = 382.7858000000001
$ws.Range("N107").Value = -5134
$ws.Range("H132").Value = 34483840
$ws.Range("I132").Value = 43478972
$ws.Range("J132").Value = 2499.8333
$ws.Range("K132").Value = 130436916
$ws.Range("L132").Value = 7499.499899999999
$ws.Range("M132").Value = -130434386
$ws.Range("N132").Value = -12559.4999
$ws.Range("H133").Value = 91787.164
$ws.Range("J133").Value = 91787.164
$ws.Range("L133").Value = 91787.164
$ws.Range("N133").Value = -101907.164
$ws.Range("H137").Value = 2478.75
$ws.Range("I137").Value = 2711.5
$ws.Range("J137").Value = 2129.625
$ws.Range("K137").Value = 8134.5
$ws.Range("L137").Value = 6388.875
$ws.Range("M137").Value = -5584.5
$ws.Range("N137").Value = -11488.875
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 3278.8
$ws.Range("I63").Value = 2497
$ws.Range("J63").Value = 3800
$ws.Range("K63").Value = 2497
$ws.Range("L63").Value = 3800
$ws.Range("M63").Value = -1811
$ws.Range("N63").Value = -5172
$ws.Range("H66").Value = 3278.8
$ws.Range("I66").Value = 2497
$ws.Range("J66").Value = 3800
$ws.Range("K66").Value = 12485
$ws.Range("L66").Value = 19000
$ws.Range("M66").Value = -9053
$ws.Range("N66").Value = -25864
$ws.Range("H74").Value = 2929.6365
$ws.Range("J74").Value = 4752.125
$ws.Range("L74").Value = 4752.125
$ws.Range("N74").Value = -6500.125
$ws.Range("H77").Value = 2929.6365
$ws.Range("J77").Value = 4752.125
$ws.Range("L77").Value = 23760.625
$ws.Range("N77").Value = -32496.625
$ws.Range("H122").Value = 4548765
$ws.Range("I122").Value = 4654399
$ws.Range("K122").Value = 13963197
$ws.Range("M122").Value = -13960747
$ws.Range("H132").Value = 5916.5835
$ws.Range("I132").Value = 3552.7727
$ws.Range("K132").Value = 10658.3181
$ws.Range("M132").Value = -8128.3181
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 76925220
$ws.Range("I20").Value = 100001900
$ws.Range("J20").Value = 2933
$ws.Range("K20").Value = 100001900
$ws.Range("L20").Value = 2933
$ws.Range("M20").Value = -100001653
$ws.Range("N20").Value = -3427
$ws.Range("H81").Value = 9795
$ws.Range("J81").Value = 9795
$ws.Range("L81").Value = 9795
$ws.Range("N81").Value = -11917
$ws.Range("H84").Value = 9795
$ws.Range("J84").Value = 9795
$ws.Range("L84").Value = 29385
$ws.Range("N84").Value = -39993
$ws.Range("H95").Value = 52277.5
$ws.Range("J95").Value = 52277.5
$ws.Range("L95").Value = 52277.5
$ws.Range("N95").Value = -57769.5
$ws.Range("H134").Value = 15152997
$ws.Range("I134").Value = 15626496
$ws.Range("K134").Value = 46879488
$ws.Range("M134").Value = -46876953
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1322.95
$ws.Range("I31").Value = 1050.4736
$ws.Range("K31").Value = 1050.4736
$ws.Range("M31").Value = -755.4736
$ws.Range("H34").Value = 1322.95
$ws.Range("I34").Value = 1050.4736
$ws.Range("K34").Value = 1050.4736
$ws.Range("M34").Value = -848.4736
$ws.Range("H62").Value = 4469.8
$ws.Range("I62").Value = 5498.5713
$ws.Range("J62").Value = 2069.3333
$ws.Range("K62").Value = 5498.5713
$ws.Range("L62").Value = 2069.3333
$ws.Range("M62").Value = -4874.5713
$ws.Range("N62").Value = -3317.3333
$ws.Range("H65").Value = 4469.8
$ws.Range("I65").Value = 5498.5713
$ws.Range("J65").Value = 2069.3333
$ws.Range("K65").Value = 27492.8565
$ws.Range("L65").Value = 10346.6665
$ws.Range("M65").Value = -24372.8565
$ws.Range("N65").Value = -16586.6665
$ws.Range("H94").Value = 1000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1000
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -1902
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4749.75
$ws.Range("I80").Value = 4898.6
$ws.Range("K80").Value = 4898.6
$ws.Range("M80").Value = -3900.6
$ws.Range("H83").Value = 4749.75
$ws.Range("I83").Value = 4898.6
$ws.Range("K83").Value = 24493
$ws.Range("M83").Value = -19501
$ws.Range("H92").Value = 108000
$ws.Range("J92").Value = 108000
$ws.Range("L92").Value = 108000
$ws.Range("N92").Value = -111744
$ws.Range("H93").Value = 44999
$ws.Range("J93").Value = 44999
$ws.Range("L93").Value = 44999
$ws.Range("N93").Value = -48743
$ws.Range("H132").Value = 1665.5238
$ws.Range("I132").Value = 1220.9445
$ws.Range("K132").Value = 3662.8335
$ws.Range("M132").Value = -1132.8335
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1108.909
$ws.Range("I46").Value = 941.5
$ws.Range("J46").Value = 1309.8
$ws.Range("K46").Value = 941.5
$ws.Range("L46").Value = 1309.8
$ws.Range("M46").Value = -753.5
$ws.Range("N46").Value = -1685.8
$ws.Range("H82").Value = 43480904
$ws.Range("I82").Value = 76924910
$ws.Range("J82").Value = 3691.6
$ws.Range("K82").Value = 76924910
$ws.Range("L82").Value = 3691.6
$ws.Range("M82").Value = -76924549
$ws.Range("N82").Value = -4413.6
$ws.Range("H85").Value = 43480904
$ws.Range("I85").Value = 76924910
$ws.Range("J85").Value = 3691.6
$ws.Range("K85").Value = 76924910
$ws.Range("L85").Value = 3691.6
$ws.Range("M85").Value = -76923662
$ws.Range("N85").Value = -6187.6
$ws.Range("H122").Value = 6699.2964
$ws.Range("I122").Value = 4173.174
$ws.Range("K122").Value = 12519.522
$ws.Range("M122").Value = -10069.522
$ws.Range("H136").Value = 2822.0557
$ws.Range("I136").Value = 2586.2666
$ws.Range("K136").Value = 7758.7998
$ws.Range("M136").Value = -5208.7998
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 31265834
$ws.Range("J62").Value = 7292.4287
$ws.Range("L62").Value = 7292.4287
$ws.Range("N62").Value = -8540.4287
$ws.Range("H65").Value = 31265834
$ws.Range("J65").Value = 7292.4287
$ws.Range("L65").Value = 36462.14350000001
$ws.Range("N65").Value = -42702.14350000001
$ws.Range("H95").Value = 25999.5
$ws.Range("J95").Value = 25999.5
$ws.Range("L95").Value = 25999.5
$ws.Range("N95").Value = -31491.5
$ws.Range("H132").Value = 2369.5
$ws.Range("I132").Value = 2369.5
$ws.Range("K132").Value = 7108.5
$ws.Range("M132").Value = -4578.5
$ws.Range("H136").Value = 13365.917
$ws.Range("I136").Value = 13365.917
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 40097.751
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -37547.751
$ws.Range("N136").ClearContents()
